$wb = $excel.ActiveWorkbook

# Map of sheet name -> row -> new F-column value ("想去人数")
$updates = @{
    "展览"   = @{ 2 = 13566; 3 = 319; 4 = 658; 5 = 226; 6 = 465; 7 = 1370; 8 = 130 }
    "全部类型" = @{ 2 = 13566; 3 = 319; 4 = 658; 5 = 226; 8 = 465; 9 = 1370; 11 = 130 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item([int]$row, 6).Value = $rows[$row]
    }
}
